$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting existing rows 6-10 down to 7-11
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the weekly entry (copy of row pattern, newer date)
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44818
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11500
$ws.Range("Q6").Value = '$/bandeja 2 kilos'
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 5750
$ws.Range("T6").Value = 2
